{"js": "// Apply the \"answers\" update: refresh the date line and all 25 three-digit\n// x one-digit multiplication answers in the table with the new values from\n// the commit. Each (old -> new) pair is unique at the moment it is applied\n// (processed in document order), so a simple search-and-replace per pair is\n// safe even though a later pair's replacement text collides with an\n// earlier pair's original text.\nconst replacements = [\n  [\"2025-09-17 Wednesday\", \"2025-09-18 Thursday\"],\n  [\"936\\u00D72=1872\", \"648\\u00D78=5184\"],\n  [\"210\\u00D78=1680\", \"625\\u00D79=5625\"],\n  [\"625\\u00D77=4375\", \"365\\u00D72=730\"],\n  [\"639\\u00D74=2556\", \"199\\u00D75=995\"],\n  [\"168\\u00D72=336\", \"384\\u00D75=1920\"],\n  [\"125\\u00D73=375\", \"919\\u00D77=6433\"],\n  [\"517\\u00D73=1551\", \"558\\u00D78=4464\"],\n  [\"331\\u00D75=1655\", \"286\\u00D75=1430\"],\n  [\"305\\u00D74=1220\", \"715\\u00D73=2145\"],\n  [\"174\\u00D76=1044\", \"572\\u00D72=1144\"],\n  [\"500\\u00D78=4000\", \"304\\u00D78=2432\"],\n  [\"670\\u00D73=2010\", \"952\\u00D72=1904\"],\n  [\"976\\u00D78=7808\", \"806\\u00D77=5642\"],\n  [\"762\\u00D76=4572\", \"733\\u00D77=5131\"],\n  [\"537\\u00D74=2148\", \"216\\u00D75=1080\"],\n  [\"857\\u00D75=4285\", \"307\\u00D77=2149\"],\n  [\"862\\u00D72=1724\", \"976\\u00D77=6832\"],\n  [\"614\\u00D78=4912\", \"294\\u00D72=588\"],\n  [\"257\\u00D77=1799\", \"228\\u00D76=1368\"],\n  [\"191\\u00D79=1719\", \"991\\u00D78=7928\"],\n  [\"421\\u00D72=842\", \"812\\u00D72=1624\"],\n  [\"265\\u00D79=2385\", \"762\\u00D76=4572\"],\n  [\"394\\u00D76=2364\", \"267\\u00D77=1869\"],\n  [\"567\\u00D73=1701\", \"445\\u00D73=1335\"],\n  [\"392\\u00D73=1176\", \"777\\u00D75=3885\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  // Replace every occurrence found for this old value (expected to be a\n  // single hit per value in this document).\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"answers\" update: refresh the date line and all 25 three-digit\n# x one-digit multiplication answers in the table with the new values from\n# the commit. Each (old -> new) pair is unique in the document at the\n# moment it is applied (processed top-to-bottom, matching document order),\n# so a straightforward Find/Replace per pair is safe even though a later\n# pair's replacement text collides with an earlier pair's original text.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-09-17 Wednesday\", \"2025-09-18 Thursday\"),\n    @(\"936\u00d72=1872\", \"648\u00d78=5184\"),\n    @(\"210\u00d78=1680\", \"625\u00d79=5625\"),\n    @(\"625\u00d77=4375\", \"365\u00d72=730\"),\n    @(\"639\u00d74=2556\", \"199\u00d75=995\"),\n    @(\"168\u00d72=336\", \"384\u00d75=1920\"),\n    @(\"125\u00d73=375\", \"919\u00d77=6433\"),\n    @(\"517\u00d73=1551\", \"558\u00d78=4464\"),\n    @(\"331\u00d75=1655\", \"286\u00d75=1430\"),\n    @(\"305\u00d74=1220\", \"715\u00d73=2145\"),\n    @(\"174\u00d76=1044\", \"572\u00d72=1144\"),\n    @(\"500\u00d78=4000\", \"304\u00d78=2432\"),\n    @(\"670\u00d73=2010\", \"952\u00d72=1904\"),\n    @(\"976\u00d78=7808\", \"806\u00d77=5642\"),\n    @(\"762\u00d76=4572\", \"733\u00d77=5131\"),\n    @(\"537\u00d74=2148\", \"216\u00d75=1080\"),\n    @(\"857\u00d75=4285\", \"307\u00d77=2149\"),\n    @(\"862\u00d72=1724\", \"976\u00d77=6832\"),\n    @(\"614\u00d78=4912\", \"294\u00d72=588\"),\n    @(\"257\u00d77=1799\", \"228\u00d76=1368\"),\n    @(\"191\u00d79=1719\", \"991\u00d78=7928\"),\n    @(\"421\u00d72=842\", \"812\u00d72=1624\"),\n    @(\"265\u00d79=2385\", \"762\u00d76=4572\"),\n    @(\"394\u00d76=2364\", \"267\u00d77=1869\"),\n    @(\"567\u00d73=1701\", \"445\u00d73=1335\"),\n    @(\"392\u00d73=1176\", \"777\u00d75=3885\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}\n\nWrite-Output \"done\"\n"}
